# ToDos.xlsx update
# - Adds 16 new Todo rows (47-62) to Sheet1 with their Status values
# - Widens column B to fit the new (longer) text
# - Updates the saved view: scroll position / active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: (row, Todo text, Status)
$newRows = @(
    @(47, "LEDs keep flashing after reset", "Open"),
    @(48, "Stabilize not listen -> listen FSM", "Closed"),
    @(49, "Stabilize startup / WiFi/ cloud connection", "Open"),
    @(50, "RateTuning", "Closed"),
    @(51, "Activate recording when incoming data", "Closed"),
    @(52, "VU meter", "Open"),
    @(53, "Add FSM diagrams", "Open"),
    @(54, "Add WiFi signal strength to cloud API and app", "Open"),
    @(55, "Used 5 band equalizer", "Open"),
    @(56, "Install stronger speaker", "Open"),
    @(57, "Only send EchoReq as keep-alive when no incoming comm.", "Open"),
    @(58, "Use exponential back-off for EchoReq", "Open"),
    @(59, "Add permission flag to cloud API to allow recording enable by remote", "Open"),
    @(60, "Add support for / experiment with ECB mode encryption", "Open"),
    @(61, "Experiment with ADPCM 16MHz mode", "Ongoing"),
    @(62, "Root cause voice_data message loss in duplex mode", "Open")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Column B needs to be a bit wider to accommodate the new/longer todo text.
$ws.Columns.Item(2).ColumnWidth = 59

# Scroll the saved view down toward the newly added rows and leave the
# selection on the last entered cell, like a user would after typing it in.
$excel.ActiveWindow.ScrollRow = 40
$null = $ws.Range("B62").Select()
